# Update the title text on slide 1 from "DAY - 2" to "DAY - 3".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$shape.TextFrame.TextRange.Text = "WEBMONK (DAY - 3) 💻"
